$wb = $excel.ActiveWorkbook

# This script applies updated FFXIV market-data values (columns H-N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the eight crafting-job
# sheets, refreshed by the scheduled market-data runner. Values are written directly
# (no formulas are used in this workbook).

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value2 = 1848.3
$ws.Cells.Item(42, 10).Value2 = 4499.5
$ws.Cells.Item(42, 12).Value2 = 13498.5
$ws.Cells.Item(42, 14).Value2 = -13958.5
$ws.Cells.Item(53, 8).Value2 = 2430.8635
$ws.Cells.Item(53, 9).Value2 = 2059
$ws.Cells.Item(53, 10).Value2 = 3227.7144
$ws.Cells.Item(53, 11).Value2 = 2059
$ws.Cells.Item(53, 12).Value2 = 3227.7144
$ws.Cells.Item(53, 13).Value2 = -1422
$ws.Cells.Item(53, 14).Value2 = -4501.7144
$ws.Cells.Item(86, 8).Value2 = 2950.3076
$ws.Cells.Item(86, 9).Value2 = 2505.4
$ws.Cells.Item(86, 10).Value2 = 3228.375
$ws.Cells.Item(86, 11).Value2 = 2505.4
$ws.Cells.Item(86, 12).Value2 = 3228.375
$ws.Cells.Item(86, 13).Value2 = -1382.4
$ws.Cells.Item(86, 14).Value2 = -5474.375
$ws.Cells.Item(89, 8).Value2 = 2950.3076
$ws.Cells.Item(89, 9).Value2 = 2505.4
$ws.Cells.Item(89, 10).Value2 = 3228.375
$ws.Cells.Item(89, 11).Value2 = 12527
$ws.Cells.Item(89, 12).Value2 = 16141.875
$ws.Cells.Item(89, 13).Value2 = -6911
$ws.Cells.Item(89, 14).Value2 = -27373.875
$ws.Cells.Item(96, 8).Value2 = 782.8461
$ws.Cells.Item(96, 9).Value2 = 862.4
$ws.Cells.Item(96, 10).Value2 = 517.6667
$ws.Cells.Item(96, 11).Value2 = 2587.2
$ws.Cells.Item(96, 12).Value2 = 1553.0001
$ws.Cells.Item(96, 13).Value2 = -1214.2
$ws.Cells.Item(96, 14).Value2 = -4299.0001
$ws.Cells.Item(97, 8).Value2 = 832
$ws.Cells.Item(97, 10).Value2 = 832
$ws.Cells.Item(97, 12).Value2 = 2496
$ws.Cells.Item(97, 14).Value2 = -3488
$ws.Cells.Item(132, 8).Value2 = 32814.53
$ws.Cells.Item(132, 9).Value2 = 36609.285
$ws.Cells.Item(132, 10).Value2 = 6251.25
$ws.Cells.Item(132, 11).Value2 = 109827.855
$ws.Cells.Item(132, 12).Value2 = 18753.75
$ws.Cells.Item(132, 13).Value2 = -107297.855
$ws.Cells.Item(132, 14).Value2 = -23813.75
$ws.Cells.Item(135, 8).Value2 = 523
$ws.Cells.Item(135, 9).Value2 = 404
$ws.Cells.Item(135, 11).Value2 = 3636
$ws.Cells.Item(135, 13).Value2 = -1101

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 3147.3125
$ws.Cells.Item(2, 9).Value2 = 1104.5
$ws.Cells.Item(2, 10).Value2 = 4373
$ws.Cells.Item(2, 11).Value2 = 1104.5
$ws.Cells.Item(2, 12).Value2 = 4373
$ws.Cells.Item(2, 13).Value2 = -991.5
$ws.Cells.Item(2, 14).Value2 = -4599
$ws.Cells.Item(32, 8).Value2 = 5528
$ws.Cells.Item(32, 9).Value2 = 5528
$ws.Cells.Item(32, 11).Value2 = 5528
$ws.Cells.Item(32, 13).Value2 = -5241
$ws.Cells.Item(74, 8).Value2 = 1234.6818
$ws.Cells.Item(74, 9).Value2 = 1063.65
$ws.Cells.Item(74, 11).Value2 = 1063.65
$ws.Cells.Item(74, 13).Value2 = -189.6500000000001
$ws.Cells.Item(77, 8).Value2 = 1234.6818
$ws.Cells.Item(77, 9).Value2 = 1063.65
$ws.Cells.Item(77, 11).Value2 = 5318.25
$ws.Cells.Item(77, 13).Value2 = -950.25
$ws.Cells.Item(116, 8).Value2 = 3147.3125
$ws.Cells.Item(116, 9).Value2 = 1104.5
$ws.Cells.Item(116, 10).Value2 = 4373
$ws.Cells.Item(116, 11).Value2 = 1104.5
$ws.Cells.Item(116, 12).Value2 = 4373
$ws.Cells.Item(116, 13).Value2 = 1189.5
$ws.Cells.Item(116, 14).Value2 = -8961
$ws.Cells.Item(132, 8).Value2 = 71430856
$ws.Cells.Item(132, 9).Value2 = 2493
$ws.Cells.Item(132, 11).Value2 = 7479
$ws.Cells.Item(132, 13).Value2 = -4949
$ws.Cells.Item(135, 8).Value2 = 34999.5
$ws.Cells.Item(135, 10).Value2 = 34999.5
$ws.Cells.Item(135, 12).Value2 = 34999.5
$ws.Cells.Item(135, 14).Value2 = -45139.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 3147.3125
$ws.Cells.Item(3, 9).Value2 = 1104.5
$ws.Cells.Item(3, 10).Value2 = 4373
$ws.Cells.Item(3, 11).Value2 = 1104.5
$ws.Cells.Item(3, 12).Value2 = 4373
$ws.Cells.Item(3, 13).Value2 = -990.5
$ws.Cells.Item(3, 14).Value2 = -4601
$ws.Cells.Item(20, 8).Value2 = 5835.2104
$ws.Cells.Item(20, 9).Value2 = 8669.362999999999
$ws.Cells.Item(20, 11).Value2 = 8669.362999999999
$ws.Cells.Item(20, 13).Value2 = -8422.362999999999
$ws.Cells.Item(86, 8).Value2 = 9881.637000000001
$ws.Cells.Item(86, 9).Value2 = 4376
$ws.Cells.Item(86, 11).Value2 = 4376
$ws.Cells.Item(86, 13).Value2 = -3253
$ws.Cells.Item(89, 8).Value2 = 9881.637000000001
$ws.Cells.Item(89, 9).Value2 = 4376
$ws.Cells.Item(89, 11).Value2 = 21880
$ws.Cells.Item(89, 13).Value2 = -16264
$ws.Cells.Item(99, 8).Value2 = 1281.6
$ws.Cells.Item(99, 9).Value2 = 1281.6
$ws.Cells.Item(99, 11).Value2 = 1281.6
$ws.Cells.Item(99, 13).Value2 = 216.4000000000001
$ws.Cells.Item(105, 8).Value2 = 1751.4615
$ws.Cells.Item(105, 9).Value2 = 1171.5
$ws.Cells.Item(105, 10).Value2 = 2679.4
$ws.Cells.Item(105, 11).Value2 = 1171.5
$ws.Cells.Item(105, 12).Value2 = 2679.4
$ws.Cells.Item(105, 13).Value2 = 575.5
$ws.Cells.Item(105, 14).Value2 = -6173.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 1582.3182
$ws.Cells.Item(31, 9).Value2 = 1759
$ws.Cells.Item(31, 10).Value2 = 1460
$ws.Cells.Item(31, 11).Value2 = 1759
$ws.Cells.Item(31, 12).Value2 = 1460
$ws.Cells.Item(31, 13).Value2 = -1464
$ws.Cells.Item(31, 14).Value2 = -2050
$ws.Cells.Item(34, 8).Value2 = 1582.3182
$ws.Cells.Item(34, 9).Value2 = 1759
$ws.Cells.Item(34, 10).Value2 = 1460
$ws.Cells.Item(34, 11).Value2 = 1759
$ws.Cells.Item(34, 12).Value2 = 1460
$ws.Cells.Item(34, 13).Value2 = -1557
$ws.Cells.Item(34, 14).Value2 = -1864
$ws.Cells.Item(58, 8).Value2 = 3466.2
$ws.Cells.Item(58, 9).Value2 = 3466.2
$ws.Cells.Item(58, 11).Value2 = 3466.2
$ws.Cells.Item(58, 13).Value2 = -3263.2
$ws.Cells.Item(105, 8).Value2 = 2815.9167
$ws.Cells.Item(105, 9).Value2 = 1691.5
$ws.Cells.Item(105, 10).Value2 = 3940.3333
$ws.Cells.Item(105, 11).Value2 = 1691.5
$ws.Cells.Item(105, 12).Value2 = 3940.3333
$ws.Cells.Item(105, 13).Value2 = 55.5
$ws.Cells.Item(105, 14).Value2 = -7434.3333
$ws.Cells.Item(132, 8).Value2 = 5485.6577
$ws.Cells.Item(132, 9).Value2 = 4377.5312
$ws.Cells.Item(132, 11).Value2 = 13132.5936
$ws.Cells.Item(132, 13).Value2 = -10602.5936
$ws.Cells.Item(136, 8).Value2 = 3466.2
$ws.Cells.Item(136, 9).Value2 = 3466.2
$ws.Cells.Item(136, 11).Value2 = 10398.6
$ws.Cells.Item(136, 13).Value2 = -7848.599999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value2 = 0
$ws.Cells.Item(22, 10).Value2 = 0
$ws.Cells.Item(22, 12).Value2 = 0
$ws.Cells.Item(27, 8).Value2 = 0
$ws.Cells.Item(27, 10).Value2 = 0
$ws.Cells.Item(27, 12).Value2 = 0
$ws.Cells.Item(122, 8).Value2 = 2436.2856
$ws.Cells.Item(122, 9).Value2 = 1916.5
$ws.Cells.Item(122, 10).Value2 = 5555
$ws.Cells.Item(122, 11).Value2 = 17248.5
$ws.Cells.Item(122, 12).Value2 = 49995
$ws.Cells.Item(122, 13).Value2 = -14798.5
$ws.Cells.Item(122, 14).Value2 = -54895
$ws.Cells.Item(132, 8).Value2 = 10450
$ws.Cells.Item(132, 9).Value2 = 2080
$ws.Cells.Item(132, 11).Value2 = 18720
$ws.Cells.Item(132, 13).Value2 = -16190
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 14).ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value2 = 5469.125
$ws.Cells.Item(122, 9).Value2 = 8130.25
$ws.Cells.Item(122, 10).Value2 = 2808
$ws.Cells.Item(122, 11).Value2 = 24390.75
$ws.Cells.Item(122, 12).Value2 = 8424
$ws.Cells.Item(122, 13).Value2 = -21940.75
$ws.Cells.Item(122, 14).Value2 = -13324

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value2 = 2056.75
$ws.Cells.Item(16, 9).Value2 = 2014.091
$ws.Cells.Item(16, 10).Value2 = 2150.6
$ws.Cells.Item(16, 11).Value2 = 2014.091
$ws.Cells.Item(16, 12).Value2 = 2150.6
$ws.Cells.Item(16, 13).Value2 = -1844.091
$ws.Cells.Item(16, 14).Value2 = -2490.6
$ws.Cells.Item(61, 8).Value2 = 4413.8335
$ws.Cells.Item(61, 9).Value2 = 3296.6
$ws.Cells.Item(61, 10).Value2 = 10000
$ws.Cells.Item(61, 11).Value2 = 3296.6
$ws.Cells.Item(61, 12).Value2 = 10000
$ws.Cells.Item(61, 13).Value2 = -3094.6
$ws.Cells.Item(61, 14).Value2 = -10404
$ws.Cells.Item(93, 8).Value2 = 14086.333
$ws.Cells.Item(93, 9).Value2 = 1808.65
$ws.Cells.Item(93, 10).Value2 = 32975.08
$ws.Cells.Item(93, 11).Value2 = 1808.65
$ws.Cells.Item(93, 12).Value2 = 32975.08
$ws.Cells.Item(93, 13).Value2 = -560.6500000000001
$ws.Cells.Item(93, 14).Value2 = -35471.08
$ws.Cells.Item(113, 8).Value2 = 4413.8335
$ws.Cells.Item(113, 9).Value2 = 3296.6
$ws.Cells.Item(113, 10).Value2 = 10000
$ws.Cells.Item(113, 11).Value2 = 3296.6
$ws.Cells.Item(113, 12).Value2 = 10000
$ws.Cells.Item(113, 13).Value2 = -1126.6
$ws.Cells.Item(113, 14).Value2 = -14340
$ws.Cells.Item(136, 8).Value2 = 66673110
$ws.Cells.Item(136, 9).Value2 = 5269.2
$ws.Cells.Item(136, 10).Value2 = 200008800
$ws.Cells.Item(136, 11).Value2 = 15807.6
$ws.Cells.Item(136, 12).Value2 = 600026400
$ws.Cells.Item(136, 13).Value2 = -13257.6
$ws.Cells.Item(136, 14).Value2 = -600031500

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value2 = 3163.25
$ws.Cells.Item(107, 9).Value2 = 1949.5
$ws.Cells.Item(107, 11).Value2 = 5848.5
$ws.Cells.Item(107, 13).Value2 = -3928.5
$ws.Cells.Item(136, 8).Value2 = 1651
$ws.Cells.Item(136, 9).Value2 = 959.75
$ws.Cells.Item(136, 11).Value2 = 2879.25
$ws.Cells.Item(136, 13).Value2 = -329.25
